$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Column width tweak: split the old B:C (min=2,max=3) merged col-width
# entry into its own per-column widths; widen column C (index 3). ---
$ws2.Columns.Item(3).ColumnWidth = 9.6666666667

# --- New header row 35 (shared strings: "dt ", "euler error", "midpoint", "RK") ---
$ws2.Range("B35").Value = "dt "
$ws2.Range("C35").Value = "euler error"
$ws2.Range("F35").Value = "midpoint"
$ws2.Range("H35").Value = "RK"

# --- Difference-quotient formulas for rows 37/38. Written before any
# PasteSpecial-based style copy below so they keep the engine's default
# (unstyled) cellXf instead of inheriting the clipboard's style. ---
$ws2.Range("D37").Formula = "=(C36-C37)/(B36-B37)"
$ws2.Range("G37").Formula = "=(F36-F37)/(B36-B37)"
$ws2.Range("I37").Formula = "=H37/(B36-B37)"
$ws2.Range("D38").Formula = "=(C37-C38)/(B37-B38)"
$ws2.Range("G38").Formula = "=(F37-F38)/(B37-B38)"
$ws2.Range("I38").Formula = "=H38/(B37-B38)"

# --- Row 36: first data row (style s=8 for B36, s=5 for C/F/H36) ---
$ws2.Range("B30").Copy()
$ws2.Range("B36").PasteSpecial(-4122)
$ws2.Range("B36").Value = [double]"2.5000000000000001E-2"

$ws2.Range("M30").Copy()
$ws2.Range("C36").PasteSpecial(-4122)
$ws2.Range("C36").Value = [double]"1.12E-4"

$ws2.Range("M30").Copy()
$ws2.Range("F36").PasteSpecial(-4122)
$ws2.Range("F36").Value = [double]"9.9999999999999995E-7"

$ws2.Range("M30").Copy()
$ws2.Range("H36").PasteSpecial(-4122)
$ws2.Range("H36").Value = [double]"1.2199999999999999E-9"

# --- Row 37 (style s=6 for B37, s=5 for C/F/H37) ---
$ws2.Range("B31").Copy()
$ws2.Range("B37").PasteSpecial(-4122)
$ws2.Range("B37").Value = [double]"1.2500000000000001E-2"

$ws2.Range("M30").Copy()
$ws2.Range("C37").PasteSpecial(-4122)
$ws2.Range("C37").Value = [double]"2.9300000000000001E-5"

$ws2.Range("M30").Copy()
$ws2.Range("F37").PasteSpecial(-4122)
$ws2.Range("F37").Value = [double]"1.2200000000000001E-7"

$ws2.Range("M30").Copy()
$ws2.Range("H37").PasteSpecial(-4122)
$ws2.Range("H37").Value = [double]"7.5600000000000003E-11"

# --- Row 38 (style s=6 for B38, s=5 for C/F/H38) ---
$ws2.Range("B31").Copy()
$ws2.Range("B38").PasteSpecial(-4122)
$ws2.Range("B38").Value = [double]"6.2500000000000003E-3"

$ws2.Range("M30").Copy()
$ws2.Range("C38").PasteSpecial(-4122)
$ws2.Range("C38").Value = [double]"1.95E-5"

$ws2.Range("M30").Copy()
$ws2.Range("F38").PasteSpecial(-4122)
$ws2.Range("F38").Value = [double]"4.06E-8"

$ws2.Range("M30").Copy()
$ws2.Range("H38").PasteSpecial(-4122)
$ws2.Range("H38").Value = [double]"7.93E-14"

$excel.CutCopyMode = 0

# --- View state: scroll sheet2 and move the selection to match the new work ---
$ws2.Activate()
$ws2.Range("F31").Select()
